$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

$wsExpo.Range("F5").Value = 507
$wsExpo.Range("F6").Value = 954
$wsExpo.Range("F7").Value = 481
$wsExpo.Range("F9").Value = 2250
$wsExpo.Range("F12").Value = 130
$wsExpo.Range("F13").Value = 1149
$wsExpo.Range("F14").Value = 193
$wsExpo.Range("F15").Value = 2263
$wsExpo.Range("F16").Value = 719
$wsExpo.Range("F17").Value = 14798
$wsExpo.Range("F19").Value = 1396
$wsExpo.Range("F22").Value = 149
$wsExpo.Range("F24").Value = 155
$wsExpo.Range("F25").Value = 112
$wsExpo.Range("F29").Value = 26
$wsExpo.Range("F31").Value = 40
$wsShow.Range("F10").Value = 22
$wsLocal.Range("F2").Value = 5750
$wsLocal.Range("F4").Value = 483
$wsAll.Range("F4").Value = 483
$wsAll.Range("F6").Value = 507
$wsAll.Range("F7").Value = 954
$wsAll.Range("F9").Value = 481
$wsAll.Range("F11").Value = 2250
$wsAll.Range("F15").Value = 130
$wsAll.Range("F17").Value = 1149
$wsAll.Range("F19").Value = 193
$wsAll.Range("F21").Value = 22
$wsAll.Range("F22").Value = 2263
$wsAll.Range("F23").Value = 719
$wsAll.Range("F27").Value = 1396
$wsAll.Range("F30").Value = 149
$wsAll.Range("F32").Value = 155
$wsAll.Range("F33").Value = 112
$wsAll.Range("F41").Value = 26
$wsAll.Range("F49").Value = 40
